$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# pre-seed with known strings
$ws.Cells.Item(1, 52).Value = "seedA1"
$ws.Cells.Item(2, 52).Value = "seedB1"
$ws.Cells.Item(1, 53).Value = "seedA2"
$ws.Cells.Item(2, 53).Value = "seedB2"

for ($col = 52; $col -le 53; $col++) {
    $cellA = $ws.Cells.Item(1, $col)
    $cellB = $ws.Cells.Item(2, $col)

    $valA = $cellA.Value()
    $valB = $cellB.Value()
    Write-Host "loaded col=$col valA=[$valA] valB=[$valB]"

    Set-TextValue $cellA $valB
    Set-TextValue $cellB $valA
}

for ($col = 52; $col -le 53; $col++) {
    $v1 = $ws.Cells.Item(1, $col).Value()
    $v2 = $ws.Cells.Item(2, $col).Value()
    Write-Host "check col=$col row1=[$v1] row2=[$v2]"
}
